$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before row 523, pushing the existing data (old rows
# 523-530) down to become rows 529-536.
$ws.Range("A523:R528").Insert()

# Common / constant values shared by all rows in this data block.
$mercadoId = 10
$mercado   = "Vega Modelo de Temuco"
$region    = "La Araucanía"
$codreg    = 9
$catId     = 100112027
$categoria = "Melón"
$unidadCom = "$/unidad"
$kgUnidad  = 1
$clasif    = "Hortaliza"

function Set-DataRow {
    param($row, $fecha, $variedad, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $origen, $precioKg)

    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $catId
    $ws.Cells.Item($row, 7).Value = $categoria
    $ws.Cells.Item($row, 8).Value = $variedad
    $ws.Cells.Item($row, 9).Value = $calidad
    $ws.Cells.Item($row, 10).Value = $volumen
    $ws.Cells.Item($row, 11).Value = $precioMin
    $ws.Cells.Item($row, 12).Value = $precioMax
    $ws.Cells.Item($row, 13).Value = $precioProm
    $ws.Cells.Item($row, 14).Value = $unidadCom
    $ws.Cells.Item($row, 15).Value = $origen
    $ws.Cells.Item($row, 16).Value = $precioKg
    $ws.Cells.Item($row, 17).Value = $kgUnidad
    $ws.Cells.Item($row, 18).Value = $clasif
}

Set-DataRow 523 44595 "Calameño" "Extra"   1500 1200 1200 1200 "Región del Maule" 1200
Set-DataRow 524 44595 "Calameño" "Primera" 5000 1000 1000 1000 "Región del Maule" 1000
Set-DataRow 525 44595 "Calameño" "Segunda" 1000 800  800  800  "Región del Maule" 800
Set-DataRow 526 44595 "Plátano"  "Primera" 500  2000 2000 2000 "Región del Maule" 2000
Set-DataRow 527 44595 "Tuna"     "Extra"   500  1200 1200 1200 "Región del Maule" 1200
Set-DataRow 528 44595 "Tuna"     "Primera" 1500 1000 1000 1000 "Región del Maule" 1000
